# "Started npc index from 0"
# The NPC sheet has, for each NPC row-block (Win/Fail/Default), a "Commands to run on NPC open"
# cell (column H on the "Win" row) containing "scriptevent rod:npcComplete N" and a
# "Commands to run on close" cell (column H on the "Fail" row) containing
# "scriptevent rod:npcReplay N". These were numbered starting at 1 (off by one from the
# npcN tag/name columns, which start at 0). This change re-numbers npc0..npc11's
# Complete/Replay commands down by one so they start at 0 (npc12 is left as-is).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("NPC_automator_input")

# Row where each NPC's "Win" row starts: npc0 -> row 2, npc1 -> row 5, ... (3 rows per npc)
for ($npc = 0; $npc -le 11; $npc++) {
    $winRow = 2 + ($npc * 3)
    $failRow = $winRow + 1

    # Clear first so the assignment below always registers as a real change, even for the
    # one npc (11) whose old text coincidentally already matches the new text (due to a
    # pre-existing authoring bug where its cells pointed at npc10's shared strings).
    $ws.Range("H$winRow").Value = ""
    $ws.Range("H$winRow").Value = "scriptevent rod:npcComplete $npc"

    $ws.Range("H$failRow").Value = ""
    $ws.Range("H$failRow").Value = "scriptevent rod:npcReplay $npc"
}

# Cosmetic: update the saved view state (scroll position / active cell) on the sheet.
$ws.Activate()
$excel.ActiveWindow.ScrollRow = 7
$excel.ActiveWindow.ScrollColumn = 7
$ws.Range("H42").Select()

$wb.Save()
